$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at the top of the data (rows 2-8), pushing the existing
# data (old row 2 = 2021-... date 44536) down to row 9, etc.
$ws.Rows("2:8").Insert()

# The newly inserted rows come in with default/no formatting; copy the
# number-format / font styling from the (now shifted) former row 2 - which
# is row 9 after the insert - onto the new rows.
$fmtSrc = $ws.Range("A9:B9")
$fmtDst = $ws.Range("A2:B8")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match the row height (18pt, same as all the other data rows) for the
# newly inserted rows.
$ws.Rows("2:8").RowHeight = 18

# Fill in the 7 new days of data (most-recent-first, matching the existing
# descending-date ordering of the sheet).
$newData = @(
    @(44543, 6419310),
    @(44542, 6418911),
    @(44541, 6417289),
    @(44540, 6413287),
    @(44539, 6409218),
    @(44538, 6405004),
    @(44537, 6400706)
)

$r = 2
foreach ($pair in $newData) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# Restore the view: selection on A8, scrolled back to the top of the sheet
# (the author had scrolled to A206 previously; after the update they're
# back near the top, selecting the newly-entered A8).
$ws.Range("A8").Select() | Out-Null
